$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.2841412015873874, 1.042070573346424, 2.569661227533701, 1.603016290476707, 1.59333094653812, 51)
    3  = @(-0.008109042937529671, 1.323519731341013, 3.714640814585274, 1.927340347366099, 1.946890523951712, 50)
    4  = @(-0.2214639157450712, 1.034360074296371, 2.089489174149197, 1.445506545868678, 1.44324361610424, 49)
    5  = @(-0.08151289187847424, 0.8998183431199257, 1.899613993327096, 1.378264848759881, 1.390412022247331, 48)
    6  = @(-0.1770401647779908, 1.055389318155083, 2.149987343866462, 1.466283514149451, 1.471292494671493, 47)
    7  = @(-0.005397054805324882, 1.29408716778638, 3.181217173515479, 1.783596695869186, 1.803297289131377, 46)
    8  = @(-0.1616785476072556, 1.097639423367275, 2.157715063143125, 1.468916288677855, 1.476489078611541, 45)
    9  = @(0.05232846318138909, 0.886449672730184, 1.305484069316138, 1.142577817619499, 1.154574473665487, 44)
    10 = @(-0.1088810105591623, 1.064859504688054, 2.120173186574139, 1.45608144915528, 1.469188947141305, 43)
    11 = @(0.1284399207036865, 1.315389267634382, 2.939163877379534, 1.71439898430311, 1.73030391165523, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
